# Apply the edit described by the diff:
#   - keep the existing "GIT HUP DENEMELERI" paragraph untouched
#   - add a new, empty paragraph after it
#   - add a new paragraph containing "Githup I ogrendim " with the
#     proofing (spell-check) marks that Word leaves around the two
#     misspelled words "Githup" and "ogrendim"

$d = $word.ActiveDocument

# Collapse a range to the very end of the document's main story so the
# new content is appended right after the last existing paragraph.
$endRange = $d.Content
$endRange.Collapse(0)

$xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml">
<pkg:xmlData>
<Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships>
</pkg:xmlData>
</pkg:part>
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p/>
<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Githup</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> I </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ogrendim</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>
</w:body></w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$endRange.InsertXML($xml)
